# The "2018" sheet tracks vacation/holiday weeks (KW01..KW53 in column A,
# taken/used flag in column B). Row 57 corresponds to KW52, which was
# marked as used (1); the commit reverts that to not-used (0).
# B2 (=B5/B3), B4 (=B3-B5) and B5 (=SUM(B6:B58)) are formulas that depend
# on B57 and will recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018")
$ws.Range("B57").Value = 0
